$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '63.538.93'
$ws.Range('E2').Value = '  -5.51%  '
$ws.Range('D3').Value = '3.279.73'
$ws.Range('E3').Value = '  -6.73%  '
$ws.Range('D4').Value = "'0.996"
$ws.Range('D4').Style = 'Normal'
$ws.Range('E4').Value = '  -0.44%  '
$ws.Range('D5').Value = "'181.34"
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -9.98%  '
$ws.Range('D6').Value = "'519.89"
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -5.85%  '
$ws.Range('D7').Value = "'0.605"
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  -0.13%  '
$ws.Range('D8').Value = '3.280.96'
$ws.Range('E8').Value = '  -6.59%  '
$ws.Range('E9').Value = '  -0.01%  '
$ws.Range('D10').Value = "'0.621"
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  -5.72%  '
$ws.Range('D11').Value = "'59.27"
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  -3.36%  '
$ws.Range('D12').Value = "'0.133"
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  -7.62%  '
$ws.Range('D13').Value = "'0.0000258"
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  -6.70%  '
$ws.Range('D14').Value = "'9.14"
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  -7.23%  '
$ws.Range('D15').Value = '3.769.46'
$ws.Range('E15').Value = '  -7.26%  '
$ws.Range('D16').Value = "'0.118"
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  -5.03%  '
$ws.Range('D17').Value = '3.257.71'
$ws.Range('E17').Value = '  -6.90%  '
$ws.Range('D18').Value = "'17.61"
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  -4.85%  '
$ws.Range('D19').Value = '63.129.03'
$ws.Range('E19').Value = '  -5.62%  '
$ws.Range('D20').Value = "'11.09"
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  -6.67%  '
$ws.Range('D21').Value = "'0.958"
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -7.58%  '
$ws.Range('D22').Value = "'372.50"
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -4.26%  '
$ws.Range('B23').Value = 'RenderToken'
$ws.Range('C23').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D23').Value = "'11.25"
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  -6.38%  '
$ws.Range('B24').Value = 'PancakeSwap'
$ws.Range('C24').Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range('D24').Value = "'3.73"
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  -7.34%  '
$ws.Range('D25').Value = "'80.53"
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  -2.42%  '
$ws.Range('D26').Value = "'3.92"
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  +5.47%  '
$ws.Range('D27').Value = "'6.07"
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -1.23%  '
$ws.Range('B28').Value = 'ImmutableX'
$ws.Range('C28').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D28').Value = "'2.66"
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  -5.26%  '
$ws.Range('B29').Value = 'InternetComputer(DFINITY)'
$ws.Range('C29').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D29').Value = "'11.54"
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  -3.70%  '
$ws.Range('D30').Value = "'8.40"
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  -5.78%  '
$ws.Range('D31').Value = "'28.75"
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  -6.61%  '
$ws.Range('D32').Value = "'6.91"
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  -5.70%  '
$ws.Range('D33').Value = "'636.98"
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  -5.43%  '
$ws.Range('D34').Value = "'11.35"
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  -3.38%  '
$ws.Range('D35').Value = "'0.106"
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  -4.18%  '
$ws.Range('D36').Value = "'58.77"
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  -7.57%  '
$ws.Range('D37').Value = "'0.403"
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  -1.97%  '
$ws.Range('E38').Value = '  +0.07%  '
$ws.Range('D39').Value = "'36.85"
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  -7.04%  '
$ws.Range('D40').Value = "'0.993"
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  -0.41%  '
$ws.Range('D41').Value = "'0.126"
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  -2.53%  '
$ws.Range('D42').Value = '2.929.04'
$ws.Range('E42').Value = '  -5.99%  '
$ws.Range('D43').Value = '0.0₃0666'
$ws.Range('E43').Value = '  -6.10%  '
$ws.Range('B44').Value = 'Fetch.AI'
$ws.Range('C44').Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range('D44').Value = "'2.47"
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  -3.38%  '
$ws.Range('B45').Value = 'ThetaToken'
$ws.Range('C45').Value = 'https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta'
$ws.Range('D45').Value = "'2.71"
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  -11.61%  '
$ws.Range('D46').Value = "'2.66"
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -0.79%  '
$ws.Range('D47').Value = "'0.0397"
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -0.40%  '
$ws.Range('D48').Value = "'2.88"
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  +9.91%  '
$ws.Range('E49').Value = '  -1.13%  '
$ws.Range('B50').Value = 'dogwifhat'
$ws.Range('C50').Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range('D50').Value = "'2.55"
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  -10.56%  '
$ws.Range('B51').Value = 'ApeXProtocol'
$ws.Range('C51').Value = 'https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex'
$ws.Range('D51').Value = "'2.97"
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  -1.72%  '
